# Update cryptos list with latest scraped price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.723.53'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.601.50'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = '''211.38'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = '''0.247'
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('D10').Value = '''19.68'
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('D12').Value = '1.826.08'
$ws.Range('D13').Value = '1.608.04'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').Value = '''65.21'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').Value = '26.693.67'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '0.0₃0745'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').Value = '''210.88'
$ws.Range('E19').Value = '  +1.20%  '
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('D21').Value = '''1.01'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '''4.31'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '''8.98'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').Value = '''143.50'
$ws.Range('E25').Value = '  -1.24%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('D30').Value = '''0.0514'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('D31').Value = '''1.15'
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('D34').Value = '1.298.47'
$ws.Range('E34').Value = '  +1.82%  '
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').Value = '''0.608'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('D38').Value = '''1.17'
$ws.Range('E38').Value = '  +21.24%  '
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('D43').Value = '''0.783'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = '''63.32'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('D45').Value = '1.736.70'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').Value = '''91.03'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '''0.102'
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.0517'
$ws.Range('E49').Value = '  +1.89%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''7.40'
$ws.Range('E51').Value = '  -0.86%  '
